$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 176, shifting existing rows 176-273 down to 177-274.
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with the new record's data.
$ws.Range("A176").Value = 1
$ws.Range("B176").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C176").Value = "Arica y Parinacota"
$ws.Range("D176").Value = 44777
$ws.Range("E176").Value = 15
$ws.Range("F176").Value = "Fruta"
$ws.Range("G176").Value = 100108
$ws.Range("H176").Value = "Tropicales y subtropicales"
$ws.Range("I176").Value = 100108006
$ws.Range("J176").Value = "Plátano"
$ws.Range("K176").Value = "Sin especificar"
$ws.Range("L176").Value = "Pintón"
$ws.Range("M176").Value = 120
$ws.Range("N176").Value = 22000
$ws.Range("O176").Value = 23000
$ws.Range("P176").Value = 22500
$ws.Range("Q176").Value = "$/caja 20 kilos"
$ws.Range("R176").Value = "Ecuador"
$ws.Range("S176").Value = 1125
$ws.Range("T176").Value = 20
